# Auto-generated script to apply value updates per commit diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 3239.476
$ws.Range("I12").Value = 2821.9333
$ws.Range("J12").Value = 4283.3335
$ws.Range("K12").Value = 2821.9333
$ws.Range("L12").Value = 4283.3335
$ws.Range("M12").Value = -2651.9333
$ws.Range("N12").Value = -4623.3335
# Row 74
$ws.Range("H74").Value = 3223.077
$ws.Range("I74").Value = 2500
$ws.Range("J74").Value = 3354.5454
$ws.Range("K74").Value = 2500
$ws.Range("L74").Value = 3354.5454
$ws.Range("M74").Value = -1564
$ws.Range("N74").Value = -5226.5454
# Row 77
$ws.Range("H77").Value = 3223.077
$ws.Range("I77").Value = 2500
$ws.Range("J77").Value = 3354.5454
$ws.Range("K77").Value = 12500
$ws.Range("L77").Value = 16772.727
$ws.Range("M77").Value = -7820
$ws.Range("N77").Value = -26132.727
# Row 107
$ws.Range("H107").Value = 446.92
$ws.Range("I107").Value = 441.5
$ws.Range("J107").Value = 486.66666
$ws.Range("K107").Value = 441.5
$ws.Range("L107").Value = 486.66666
$ws.Range("M107").Value = 1478.5
$ws.Range("N107").Value = -4326.66666
# Row 116
$ws.Range("H116").Value = 7845.125
$ws.Range("I116").Value = 1192
$ws.Range("K116").Value = 1192
$ws.Range("M116").Value = 2250
# Row 125
$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 9000
$ws.Range("M125").Value = -6540
# Row 129
$ws.Range("H129").Value = 874.88464
$ws.Range("I129").Value = 397.6
$ws.Range("J129").Value = 988.5238000000001
$ws.Range("K129").Value = 1192.8
$ws.Range("L129").Value = 2965.5714
$ws.Range("M129").Value = 3807.2
$ws.Range("N129").Value = -12965.5714
# Row 130
$ws.Range("H130").Value = 32500
$ws.Range("J130").Value = 32500
$ws.Range("L130").Value = 32500
$ws.Range("N130").Value = -42540

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2206.55
$ws.Range("I2").Value = 2116.5
$ws.Range("J2").Value = 2416.6667
$ws.Range("K2").Value = 2116.5
$ws.Range("L2").Value = 2416.6667
$ws.Range("M2").Value = -2003.5
$ws.Range("N2").Value = -2642.6667
# Row 46
$ws.Range("H46").Value = 82307.84
$ws.Range("I46").Value = 6158.3335
$ws.Range("J46").Value = 147578.86
$ws.Range("K46").Value = 6158.3335
$ws.Range("L46").Value = 147578.86
$ws.Range("M46").Value = -5839.3335
$ws.Range("N46").Value = -148216.86
# Row 74
$ws.Range("H74").Value = 126439
$ws.Range("I74").Value = 144244.58
$ws.Range("J74").Value = 1800
$ws.Range("K74").Value = 144244.58
$ws.Range("L74").Value = 1800
$ws.Range("M74").Value = -143370.58
$ws.Range("N74").Value = -3548
# Row 77
$ws.Range("H77").Value = 126439
$ws.Range("I77").Value = 144244.58
$ws.Range("J77").Value = 1800
$ws.Range("K77").Value = 721222.8999999999
$ws.Range("L77").Value = 9000
$ws.Range("M77").Value = -716854.8999999999
$ws.Range("N77").Value = -17736
# Row 110
$ws.Range("H110").Value = 1422
$ws.Range("I110").Value = 1542.3077
$ws.Range("J110").Value = 640
$ws.Range("K110").Value = 1542.3077
$ws.Range("L110").Value = 640
$ws.Range("M110").Value = 502.6922999999999
$ws.Range("N110").Value = -4730
# Row 116
$ws.Range("H116").Value = 2206.55
$ws.Range("I116").Value = 2116.5
$ws.Range("J116").Value = 2416.6667
$ws.Range("K116").Value = 2116.5
$ws.Range("L116").Value = 2416.6667
$ws.Range("M116").Value = 177.5
$ws.Range("N116").Value = -7004.6667

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2206.55
$ws.Range("I3").Value = 2116.5
$ws.Range("J3").Value = 2416.6667
$ws.Range("K3").Value = 2116.5
$ws.Range("L3").Value = 2416.6667
$ws.Range("M3").Value = -2002.5
$ws.Range("N3").Value = -2644.6667
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 39
$ws.Range("H39").Value = 2051
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
# Row 49
$ws.Range("H49").Value = 2051
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
# Row 107
$ws.Range("H107").Value = 1929.7142
$ws.Range("I107").Value = 473.33334
$ws.Range("J107").Value = 3022
$ws.Range("K107").Value = 473.33334
$ws.Range("L107").Value = 3022
$ws.Range("M107").Value = 1446.66666
$ws.Range("N107").Value = -6862
# Row 132
$ws.Range("H132").Value = 3268.25
$ws.Range("I132").Value = 2690.0625
$ws.Range("K132").Value = 8070.1875
$ws.Range("M132").Value = -5540.1875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 53292.125
$ws.Range("J37").Value = 53292.125
$ws.Range("L37").Value = 159876.375
$ws.Range("N37").Value = -160100.375
# Row 69
$ws.Range("H69").Value = 1080
$ws.Range("I69").Value = 542.8570999999999
$ws.Range("J69").Value = 2333.3333
$ws.Range("K69").Value = 1628.5713
$ws.Range("L69").Value = 6999.999899999999
$ws.Range("M69").Value = -817.5712999999998
$ws.Range("N69").Value = -8621.999899999999
# Row 72
$ws.Range("H72").Value = 1080
$ws.Range("I72").Value = 542.8570999999999
$ws.Range("J72").Value = 2333.3333
$ws.Range("K72").Value = 4885.7139
$ws.Range("L72").Value = 20999.9997
$ws.Range("M72").Value = -829.7138999999997
$ws.Range("N72").Value = -29111.9997
# Row 131
$ws.Range("H131").Value = 867.62
$ws.Range("J131").Value = 897.36664
$ws.Range("L131").Value = 2692.09992
$ws.Range("N131").Value = -12772.09992
# Row 139
$ws.Range("H139").Value = 1785.6154
$ws.Range("I139").Value = 782.8570999999999
$ws.Range("J139").Value = 2955.5
$ws.Range("K139").Value = 2348.5713
$ws.Range("L139").Value = 8866.5
$ws.Range("M139").Value = 2791.4287
$ws.Range("N139").Value = -19146.5
# Row 141
$ws.Range("H141").Value = 6620
$ws.Range("I141").Value = 3240
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 9720
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -4540
$ws.Range("N141").Value = -40360

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
# Row 130
$ws.Range("H130").Value = 22500
$ws.Range("J130").Value = 22500
$ws.Range("L130").Value = 22500
$ws.Range("N130").Value = -32540
# Row 132
$ws.Range("H132").Value = 5261.6875
$ws.Range("I132").Value = 6108
$ws.Range("J132").Value = 4415.375
$ws.Range("K132").Value = 18324
$ws.Range("L132").Value = 13246.125
$ws.Range("M132").Value = -15794
$ws.Range("N132").Value = -18306.125
# Row 133
$ws.Range("H133").Value = 31214
$ws.Range("J133").Value = 31214
$ws.Range("L133").Value = 31214
$ws.Range("N133").Value = -36274

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1064.4445
$ws.Range("I81").Value = 1064.4445
$ws.Range("K81").Value = 2128.889
$ws.Range("M81").Value = -1067.889
# Row 84
$ws.Range("H84").Value = 1064.4445
$ws.Range("I84").Value = 1064.4445
$ws.Range("K84").Value = 10644.445
$ws.Range("M84").Value = -5340.445
